# Actualización automática del tracker
# Append two new result rows (7 and 8) to the results tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
$ws.Cells.Item(7, 1).Value = 14552618
# "fecha" holds a plain text date string (e.g. "2025-09-06"), not a real
# Excel date. Force text formatting before the write so it is not
# auto-converted into a date serial, then drop the temporary format so the
# cell is left unstyled, matching the other data rows.
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "2025-09-06"
$ws.Cells.Item(7, 2).ClearFormats()
$ws.Cells.Item(7, 3).Value = "Santiago Rodriguez Taverna"
$ws.Cells.Item(7, 4).Value = "Andrew Paulson"
$ws.Cells.Item(7, 5).Value = "Gana Andrew Paulson"
$ws.Cells.Item(7, 6).Value = 1.91
# "resultado" / "profit" are still pending (no match result yet), same as
# the placeholder cells already present on row 6.
$ws.Cells.Item(7, 7).Value = ""
$ws.Cells.Item(7, 8).Value = ""

# --- Row 8 ---
$ws.Cells.Item(8, 1).Value = 14552663
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "2025-09-06"
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 3).Value = "Alex Molcan"
$ws.Cells.Item(8, 4).Value = "Stefanos Sakellaridis"
$ws.Cells.Item(8, 5).Value = "Gana Stefanos Sakellaridis"
$ws.Cells.Item(8, 6).Value = 2.75
$ws.Cells.Item(8, 7).Value = ""
$ws.Cells.Item(8, 8).Value = ""
